# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from serial date 45188 (2023-09-19) to 45189 (2023-09-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
